$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2,4)
$c.NumberFormat = "@"
$c.Value = "33.893.74"
$c.Style = "Normal"
$ws.Cells.Item(2,5).Value = "  +7.21%  "

$c = $ws.Cells.Item(3,4)
$c.NumberFormat = "@"
$c.Value = "1.779.01"
$c.Style = "Normal"
$ws.Cells.Item(3,5).Value = "  +3.65%  "

$c = $ws.Cells.Item(4,4)
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Cells.Item(4,5).Value = "  -0.20%  "

$c = $ws.Cells.Item(5,4)
$c.NumberFormat = "@"
$c.Value = "224.68"
$c.Style = "Normal"
$ws.Cells.Item(5,5).Value = "  +0.37%  "

$c = $ws.Cells.Item(6,4)
$c.NumberFormat = "@"
$c.Value = "0.559"
$c.Style = "Normal"
$ws.Cells.Item(6,5).Value = "  +3.73%  "

$c = $ws.Cells.Item(7,4)
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Cells.Item(7,5).Value = "  -0.05%  "

$c = $ws.Cells.Item(8,4)
$c.NumberFormat = "@"
$c.Value = "30.14"
$c.Style = "Normal"
$ws.Cells.Item(8,5).Value = "  +0.36%  "

$c = $ws.Cells.Item(9,4)
$c.NumberFormat = "@"
$c.Value = "46.59"
$c.Style = "Normal"
$ws.Cells.Item(9,5).Value = "  +3.81%  "

$c = $ws.Cells.Item(10,4)
$c.NumberFormat = "@"
$c.Value = "0.278"
$c.Style = "Normal"
$ws.Cells.Item(10,5).Value = "  +2.87%  "

$c = $ws.Cells.Item(11,4)
$c.NumberFormat = "@"
$c.Value = "0.0663"
$c.Style = "Normal"
$ws.Cells.Item(11,5).Value = "  +1.05%  "

$ws.Cells.Item(12,5).Value = "  +1.05%  "

$c = $ws.Cells.Item(13,4)
$c.NumberFormat = "@"
$c.Value = "2.033.81"
$c.Style = "Normal"
$ws.Cells.Item(13,5).Value = "  +3.52%  "

$c = $ws.Cells.Item(14,4)
$c.NumberFormat = "@"
$c.Value = "1.781.07"
$c.Style = "Normal"
$ws.Cells.Item(14,5).Value = "  +3.66%  "

$c = $ws.Cells.Item(15,4)
$c.NumberFormat = "@"
$c.Value = "0.623"
$c.Style = "Normal"
$ws.Cells.Item(15,5).Value = "  +1.20%  "

$c = $ws.Cells.Item(16,4)
$c.NumberFormat = "@"
$c.Value = "33.875.80"
$c.Style = "Normal"
$ws.Cells.Item(16,5).Value = "  +7.18%  "

$c = $ws.Cells.Item(17,4)
$c.NumberFormat = "@"
$c.Value = "10.11"
$c.Style = "Normal"
$ws.Cells.Item(17,5).Value = "  -0.62%  "

$c = $ws.Cells.Item(18,4)
$c.NumberFormat = "@"
$c.Value = "4.17"
$c.Style = "Normal"
$ws.Cells.Item(18,5).Value = "  -0.41%  "

$c = $ws.Cells.Item(19,4)
$c.NumberFormat = "@"
$c.Value = "68.42"
$c.Style = "Normal"
$ws.Cells.Item(19,5).Value = "  +1.31%  "

$c = $ws.Cells.Item(20,4)
$c.NumberFormat = "@"
$c.Value = "249.93"
$c.Style = "Normal"
$ws.Cells.Item(20,5).Value = "  -0.66%  "

$c = $ws.Cells.Item(21,4)
$c.NumberFormat = "@"
$c.Value = "0.0₃0738"
$c.Style = "Normal"
$ws.Cells.Item(21,5).Value = "  +1.38%  "

$c = $ws.Cells.Item(22,4)
$c.NumberFormat = "@"
$c.Value = "0.997"
$c.Style = "Normal"
$ws.Cells.Item(22,5).Value = "  -0.17%  "

$c = $ws.Cells.Item(23,4)
$c.NumberFormat = "@"
$c.Value = "10.29"
$c.Style = "Normal"
$ws.Cells.Item(23,5).Value = "  +1.30%  "

$c = $ws.Cells.Item(24,4)
$c.NumberFormat = "@"
$c.Value = "4.15"
$c.Style = "Normal"
$ws.Cells.Item(24,5).Value = "  -2.64%  "

$ws.Cells.Item(25,5).Value = "  -1.86%  "

$c = $ws.Cells.Item(26,4)
$c.NumberFormat = "@"
$c.Value = "158.43"
$c.Style = "Normal"
$ws.Cells.Item(26,5).Value = "  -0.46%  "

$c = $ws.Cells.Item(27,4)
$c.NumberFormat = "@"
$c.Value = "16.44"
$c.Style = "Normal"
$ws.Cells.Item(27,5).Value = "  +1.76%  "

$ws.Cells.Item(28,5).Value = "  +0.55%  "

$c = $ws.Cells.Item(29,4)
$c.NumberFormat = "@"
$c.Value = "6.95"
$c.Style = "Normal"
$ws.Cells.Item(29,5).Value = "  +1.68%  "

$c = $ws.Cells.Item(30,4)
$c.NumberFormat = "@"
$c.Value = "0.996"
$c.Style = "Normal"
$ws.Cells.Item(30,5).Value = "  -0.38%  "

$c = $ws.Cells.Item(31,4)
$c.NumberFormat = "@"
$c.Value = "3.81"
$c.Style = "Normal"
$ws.Cells.Item(31,5).Value = "  -1.97%  "

$c = $ws.Cells.Item(32,4)
$c.NumberFormat = "@"
$c.Value = "0.0515"
$c.Style = "Normal"
$ws.Cells.Item(32,5).Value = "  +2.23%  "

$ws.Cells.Item(33,5).Value = "  +2.07%  "

$c = $ws.Cells.Item(34,4)
$c.NumberFormat = "@"
$c.Value = "3.57"
$c.Style = "Normal"
$ws.Cells.Item(34,5).Value = "  +3.92%  "

$c = $ws.Cells.Item(35,4)
$c.NumberFormat = "@"
$c.Value = "1.83"
$c.Style = "Normal"
$ws.Cells.Item(35,5).Value = "  +4.57%  "

$c = $ws.Cells.Item(36,4)
$c.NumberFormat = "@"
$c.Value = "1.489.53"
$c.Style = "Normal"
$ws.Cells.Item(36,5).Value = "  -2.62%  "

$c = $ws.Cells.Item(37,4)
$c.NumberFormat = "@"
$c.Value = "1.07"
$c.Style = "Normal"
$ws.Cells.Item(37,5).Value = "  +1.78%  "

$c = $ws.Cells.Item(38,4)
$c.NumberFormat = "@"
$c.Value = "0.632"
$c.Style = "Normal"
$ws.Cells.Item(38,5).Value = "  +3.43%  "

$ws.Cells.Item(39,2).Value = "VeChain"
$ws.Cells.Item(39,3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Cells.Item(39,4)
$c.NumberFormat = "@"
$c.Value = "0.0185"
$c.Style = "Normal"
$ws.Cells.Item(39,5).Value = "  +1.61%  "

$ws.Cells.Item(40,2).Value = "Aave"
$ws.Cells.Item(40,3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Cells.Item(40,4)
$c.NumberFormat = "@"
$c.Value = "83.10"
$c.Style = "Normal"
$ws.Cells.Item(40,5).Value = "  -0.13%  "

$c = $ws.Cells.Item(41,4)
$c.NumberFormat = "@"
$c.Value = "2.34"
$c.Style = "Normal"
$ws.Cells.Item(41,5).Value = "  +1.40%  "

$c = $ws.Cells.Item(42,4)
$c.NumberFormat = "@"
$c.Value = "2.71"
$c.Style = "Normal"
$ws.Cells.Item(42,5).Value = "  -0.41%  "

$c = $ws.Cells.Item(43,4)
$c.NumberFormat = "@"
$c.Value = "0.888"
$c.Style = "Normal"
$ws.Cells.Item(43,5).Value = "  +3.80%  "

$c = $ws.Cells.Item(44,4)
$c.NumberFormat = "@"
$c.Value = "2.06"
$c.Style = "Normal"
$ws.Cells.Item(44,5).Value = "  +1.10%  "

$ws.Cells.Item(45,5).Value = "  +2.34%  "

$ws.Cells.Item(46,5).Value = "  +4.09%  "

$c = $ws.Cells.Item(47,4)
$c.NumberFormat = "@"
$c.Value = "1.927.99"
$c.Style = "Normal"
$ws.Cells.Item(47,5).Value = "  +3.66%  "

$ws.Cells.Item(48,2).Value = "FraxShare"
$ws.Cells.Item(48,3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Cells.Item(48,4)
$c.NumberFormat = "@"
$c.Value = "5.70"
$c.Style = "Normal"
$ws.Cells.Item(48,5).Value = "  +1.60%  "

$ws.Cells.Item(49,2).Value = "PaxDollar"
$ws.Cells.Item(49,3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$c = $ws.Cells.Item(49,4)
$c.NumberFormat = "@"
$c.Value = "0.997"
$c.Style = "Normal"
$ws.Cells.Item(49,5).Value = "  -0.27%  "

$c = $ws.Cells.Item(50,4)
$c.NumberFormat = "@"
$c.Value = "11.81"
$c.Style = "Normal"
$ws.Cells.Item(50,5).Value = "  +14.25%  "

$c = $ws.Cells.Item(51,4)
$c.NumberFormat = "@"
$c.Value = "51.36"
$c.Style = "Normal"
$ws.Cells.Item(51,5).Value = "  -2.88%  "
